# ---------------------------------------------------------------------------
# Edit script: reproduces the "crispianz" commit on crispian.docx
#   1. Append a red "(This is a change - Version for main branch)" note
#      (split into 3 runs) after the first paragraph's existing text, which
#      itself gains two trailing spaces.
#   2. Append a new empty paragraph (style "larger") after the final
#      paragraph of the Crispin's day speech, right before the section break.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# --- 1. First paragraph: add trailing spaces + red annotation -------------
$p1 = $d.Paragraphs(1)
$r1 = $p1.Range
# End of the paragraph's text range sits just before the paragraph mark.
$insertPoint = $r1.End - 1

$spaces = $d.Range($insertPoint, $insertPoint)
$spaces.InsertAfter("  ")
$insertPoint = $insertPoint + 2

$dash = [char]0x2013

$seg1Text = "(This is a change " + $dash + " Ve"
$seg1 = $d.Range($insertPoint, $insertPoint)
$seg1.InsertAfter($seg1Text)
$seg1 = $d.Range($insertPoint, $insertPoint + $seg1Text.Length)
$seg1.Font.Color = 255
$insertPoint = $insertPoint + $seg1Text.Length

$seg2Text = "rsion for main branch"
$seg2 = $d.Range($insertPoint, $insertPoint)
$seg2.InsertAfter($seg2Text)
$seg2 = $d.Range($insertPoint, $insertPoint + $seg2Text.Length)
$seg2.Font.Color = 255
$insertPoint = $insertPoint + $seg2Text.Length

$seg3Text = ")"
$seg3 = $d.Range($insertPoint, $insertPoint)
$seg3.InsertAfter($seg3Text)
$seg3 = $d.Range($insertPoint, $insertPoint + $seg3Text.Length)
$seg3.Font.Color = 255
$insertPoint = $insertPoint + $seg3Text.Length

# --- 2. Append new empty paragraph (style "larger") at very end -----------
$lastParaIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastParaIndex)
$endOfDoc = $lastPara.Range.End

$endRange = $d.Range($endOfDoc, $endOfDoc)
$endRange.InsertParagraphAfter()

$newParaIndex = $d.Paragraphs.Count
$newPara = $d.Paragraphs($newParaIndex)
$newPara.Style = "larger"
$newPara.Format.Shading.BackgroundPatternColor = 0xFFFFFF
$newPara.Format.SpaceBefore = 0
$newPara.Format.SpaceBeforeAuto = $false
$newPara.Format.SpaceAfter = 7.5
$newPara.Format.SpaceAfterAuto = $false

Write-Output "done"
